$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.507.20"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "1.871.69"
$ws.Range("E3").Value = "  -0.14%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.009"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -1.54%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.21"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.68%  "
$ws.Range("E6").Value = "  -1.33%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5081"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3897"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.61%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08359"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.28%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.106"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.84"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.211"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.78%  "
$ws.Range("D13").Value = "1.869.69"
$ws.Range("E13").Value = "  +3.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.41"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.270"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.008"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.56%  "
$ws.Range("E17").Value = "  -0.86%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.16"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.29%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06736"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.44%  "
$ws.Range("E20").Value = "  -0.29%  "
$ws.Range("E21").Value = "  -1.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.921"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.11%  "
$ws.Range("D23").Value = "28.541.69"
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.09"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.54%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.210"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.68%  "
$ws.Range("D26").Value = "2.083.64"
$ws.Range("E26").Value = "  +3.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "158.59"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.94%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.60"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.426"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.22"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1038"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.51%  "
$ws.Range("E32").Value = "  +0.63%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.731"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.69%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.614"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.45%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02458"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.87%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06593"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.929"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.59%  "
$ws.Range("E38").Value = "  -1.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.035"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.181"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.43%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.240"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.88%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6374"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.14%  "
$ws.Range("E43").Value = "  -0.57%  "
$ws.Range("E44").Value = "  -1.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6001"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.79%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.07"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.19%  "
$ws.Range("E47").Value = "  -0.84%  "
$ws.Range("E48").Value = "  +0.32%  "
$ws.Range("E49").Value = "  +0.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "122.56"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.59%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06809"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.83%  "
